{"js": "// Replace the lone \"start_ack\" mention in the documentation body with\n// \"ready\" (the handshake signal was renamed), per the commit message\n// \"Removed start_ack from documentation\".\nconst results = context.document.body.search(\"start_ack\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"ready\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the lone \"start_ack\" mention in the documentation body with\n# \"ready\" (the handshake signal was renamed), per the commit message\n# \"Removed start_ack from documentation\".\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"start_ack\"\n$find.Replacement.Text = \"ready\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Wrap = 1\n\n$find.Execute(\n    [ref]\"start_ack\",\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]\"ready\",\n    [ref]2\n)\n"}
